$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 3 (pushes the existing rows 3..63 down to 4..64)
$ws.Rows(3).Insert()

# 2. Populate the new row 3 with the "Tena Silhouette waschbarer Inko-Slip M" product.
#    Columns that hold numeric-looking text in the original data (id, price fields,
#    quantity, priceContextAmount, etc.) must stay text, so force a text number
#    format before assigning them.
$textCols = @("A3","D3","H3","I3","K3","L3")
foreach ($addr in $textCols) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A3").Value = "6866483"
$ws.Range("B3").Value = "Tena Silhouette waschbarer Inko-Slip M"
$ws.Range("C3").Value = "/de/haushalt-tier/bekleidung/socken-unterwaesche/unterwaesche/tena-silhouette-waschbarer-inko-slip-m/p/6866483"
$ws.Range("D3").Value = "1ST"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = "Tena"
$ws.Range("H3").Value = "34.95"
$ws.Range("I3").Value = "34.95/1ST"
$ws.Range("J3").Value = "Preis pro 1 Stück"
$ws.Range("K3").Value = "34.95"
$ws.Range("L3").Value = "1ST"
$ws.Range("M3").Value = "['haushalt-tier', 'bekleidung', 'socken-unterwaesche', 'unterwaesche']"
$ws.Range("N3").Value = "Tena Silhouette waschbarer Inko-Slip M 34.95 Schweizer Franken"
$ws.Range("O3").Value = "2022-09-13 21:01:11"

# 3. Remove what is now the trailing row (old last row, id 6866484, shifted to row 64).
$ws.Rows(64).Delete()

# 4. Refresh the timestamp column for every remaining data row (2..63) to the new value.
for ($r = 2; $r -le 63; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-09-13 21:01:11"
}
